$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Round row 5 (columns B..AH) to 2 decimal places, matching Excel's ROUND()
# round-half-away-from-zero behaviour.
$rng = $ws.Range("B5:AH5")
foreach ($cell in $rng.Cells) {
    $v = $cell.Value2
    if ($v -ne $null) {
        $cell.Value = $excel.WorksheetFunction.Round([double]$v, 2)
    }
}

# Delete row 6 entirely (shifts rows up).
$ws.Rows.Item(6).Delete()
